# Week 15 simulations added.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet - append new per-play yardage samples to the 4 running logs
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value2 = $ydsWs.Range("B2").Value2 + " 7 8 4 4 2 0 7 -1 2 7 5 1 7 5 7 4 9 5 -3 13 8 5 8 6 8 7"
$ydsWs.Range("B3").Value2 = $ydsWs.Range("B3").Value2 + " 5 15 6 60 3 7 -2 1 17 3 9 6 18 10 1 8 24"
$ydsWs.Range("C2").Value2 = $ydsWs.Range("C2").Value2 + " 0 14 9 5 5 9 1 1 0 2 17 2 3 3 12 4 -1 6 4 8 2 4 2 6 1 3 2 10 9 0 3 2 0 2 4"
$ydsWs.Range("C3").Value2 = $ydsWs.Range("C3").Value2 + " 15 7 5 11 9 -2 15 12 10 25 7 2 7 5 17 59 7 36 5 11 8 3 1"

# ---------------------------------------------------------------------
# OFF sheet - running totals for offensive down/distance situations
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value2 = 133
$offWs.Range("D2").Value2 = 6
$offWs.Range("F2").Value2 = 59
$offWs.Range("G2").Value2 = 37
$offWs.Range("I2").Value2 = 6
$offWs.Range("J2").Value2 = 21
$offWs.Range("N2").Value2 = 14
$offWs.Range("O2").Value2 = 20

$offWs.Range("C3").Value2 = 173
$offWs.Range("E3").Value2 = 22
$offWs.Range("F3").Value2 = 88
$offWs.Range("G3").Value2 = 23
$offWs.Range("H3").Value2 = 22
$offWs.Range("I3").Value2 = 61
$offWs.Range("J3").Value2 = 42
$offWs.Range("L3").Value2 = 261
$offWs.Range("M3").Value2 = 156
$offWs.Range("Q3").Value2 = 424

# ---------------------------------------------------------------------
# DEF sheet - running totals for defensive down/distance situations
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value2 = 174
$defWs.Range("D2").Value2 = 9
$defWs.Range("F2").Value2 = 49
$defWs.Range("G2").Value2 = 62
$defWs.Range("I2").Value2 = 5
$defWs.Range("J2").Value2 = 23
$defWs.Range("N2").Value2 = 10

$defWs.Range("C3").Value2 = 152
$defWs.Range("D3").Value2 = 1
$defWs.Range("E3").Value2 = 22
$defWs.Range("F3").Value2 = 92
$defWs.Range("G3").Value2 = 33
$defWs.Range("H3").Value2 = 12
$defWs.Range("I3").Value2 = 56
$defWs.Range("J3").Value2 = 50
$defWs.Range("L3").Value2 = 272
$defWs.Range("M3").Value2 = 190
$defWs.Range("Q3").Value2 = 496

# ---------------------------------------------------------------------
# ST sheet - special teams totals + running logs
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value2 = 55
$stWs.Range("D2").Value2 = 54
$stWs.Range("F2").Value2 = 359
$stWs.Range("G2").Value2 = 342
$stWs.Range("H2").Value2 = 5
$stWs.Range("I2").Value2 = 4

$stWs.Range("B4").Value2 = $stWs.Range("B4").Value2 + " 58 71"
$stWs.Range("B5").Value2 = $stWs.Range("B5").Value2 + " 0 29"
$stWs.Range("B6").Value2 = $stWs.Range("B6").Value2 + " 17 17"
$stWs.Range("D3").Value2 = $stWs.Range("D3").Value2 + " 50 46 51 63 18"
$stWs.Range("D4").Value2 = $stWs.Range("D4").Value2 + " 0 0 0 0 0"
$stWs.Range("D5").Value2 = $stWs.Range("D5").Value2 + " 8 0 0"

# ---------------------------------------------------------------------
# TURNS sheet - turnovers
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value2 = 8
$turnsWs.Range("D3").Value2 = 7

# ---------------------------------------------------------------------
# PEN sheet - penalties
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value2 = 13
